$wb = $excel.ActiveWorkbook

# Updates to "想去人数" (column F) on sheets "展览" and "全部类型"
$updates = @{
    2  = 284
    4  = 3493
    5  = 2178
    8  = 61
    9  = 48
    10 = 1287
    12 = 1712
    13 = 125
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates[$row]
}

$updates2 = @{
    2  = 284
    4  = 3493
    5  = 2178
    9  = 61
    10 = 48
    13 = 1287
    15 = 1712
    16 = 125
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates2.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates2[$row]
}
